# Convention change to support multi-axle vehicles:
# Add a new "Truck_Amandla" sheet, duplicated from "Bus_Makhulu",
# with updated x/y/z offset values and relabeled header cell.

$wb = $excel.ActiveWorkbook

# Duplicate the Bus_Makhulu sheet and place the copy at the end of the tab strip.
$srcWs = $wb.Worksheets.Item("Bus_Makhulu")
$srcWs.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs.Name = "Truck_Amandla"

# Relabel the class/title cell (H3) to the new sheet's own name.
$newWs.Range("H3").Value = "Truck_Amandla"

# Update the x / y / z offset values on row 5 for the new vehicle class.
$newWs.Range("F5").Value = -1.3231
$newWs.Range("G5").Value = 0.558013
$newWs.Range("H5").Value = 2.3924

# Force a distinct conditional-format style entry for the new sheet
# (keeps each tab's "class" highlight rule backed by its own dxf record).
$cfRange = $newWs.Range("A4:B4")
$cfRange.FormatConditions.Delete()
$newRule = $cfRange.FormatConditions.Add(1, 3, '"class"')
$newRule.Interior.Color = 13431551

# Make the new sheet the active tab with the same selection it was saved with.
$newWs.Activate() | Out-Null
$newWs.Range("N10").Select() | Out-Null
